$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Flashcards" column data (E:G) appended to the existing Grade/Subjects/Lesson/Con table.
$ws.Range("E1").Value = "Quiz"
$ws.Range("F1").Value = "Worksheet"
$ws.Range("G1").Value = "Flashcards"

$ws.Range("E2").Value = "dsd"

$ws.Range("E3").Value = "s"

$ws.Range("F4").Value = "sdds"

$ws.Range("G5").Value = "dsds"

$ws.Range("G6").Value = "dsds"

$ws.Range("E8").Value = "dsd"
$ws.Range("F8").Value = "dsdsf"
$ws.Range("G8").Value = "fdfdfd"

# Widen column D to fit its contents (Excel snaps ColumnWidth to whole
# pixels, so 28.8 is the input that lands closest to the saved 29.6640625
# "characters" width once re-serialised).
$ws.Columns.Item(4).ColumnWidth = 28.8

# Match the saved selection state.
$ws.Range("G8").Select()
